# "Generate Report for Handoff"
# Swap the handed-off source file's GUID-named identity (and its refreshed
# handoff timestamps) throughout the localization-status report:
#   cc10638c-e9a9-4176-86b5-f50732c98327  ->  9c5658cd-12a0-4987-971b-c85449672038
# and update the three "Latest * Datetime" stamps that the handoff run
# produced.

$wb = $excel.ActiveWorkbook

$oldGuid = "cc10638c-e9a9-4176-86b5-f50732c98327"
$newGuid = "9c5658cd-12a0-4987-971b-c85449672038"

$oldZhHash = "64b3c01a956b6566e2c3973757a537b216ac1002"
$newZhHash = "6c2a7895ea9fe90145e9ea84973028a5db512e4a"

# The external hyperlink target (unchanged by this edit - only the visible
# display text moves to the new file name) - identical on all 3 sheets.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d77e3f5888053352ff2822dedf4c6ab623c2efee/e2e/$oldGuid.md"

function Set-HyperlinkDisplay($range, [string]$address, [string]$display) {
    # Re-stamp the hyperlink on $range with new display text while keeping
    # the same target address. Deleting the whole Hyperlinks collection for
    # the range first (rather than editing TextToDisplay in place) avoids
    # leaving a stale duplicate <hyperlink> entry behind.
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $address, "", "", $display) | Out-Null
    # Hyperlinks.Add re-stamps the range with the workbook's built-in
    # "Hyperlink" cell style; put back the original look (underlined,
    # cornflower blue 0x6495ED) that the report's hand-rolled HyperLink
    # style used, so only the text content moves - not the formatting.
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

Set-HyperlinkDisplay $wsOverview.Range("B2") $linkAddress "e2e\$newGuid.md"

$wsOverview.Range("G2").Value = "2016-10-13 13:35:00"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HyperlinkDisplay $wsZh.Range("A2") $linkAddress "$newGuid.md"

$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-10-13 13:34:50"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

Set-HyperlinkDisplay $wsDe.Range("A2") $linkAddress "$newGuid.md"

$wsDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-10-13 13:35:00"
